# Crowdin sync update for ja/Backup.xlsx
#
# The Japanese locale sheet previously carried three columns per entry:
#   A = translation key, B = Japanese text, E = English reference text.
# The Crowdin export for this locale drops the English reference column
# (column E) for the Coroner / dead-body related keys, keeping only the
# key (A) and the Japanese translation (B).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Rows that contained a "reference" English string in column E that must
# be removed now that only key/Japanese-text pairs are kept.
$rowsToClear = @(26, 27, 28, 29, 30, 31, 33, 34, 35, 36, 38, 39, 40, 41, 43)

foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 5).ClearContents()
}

# Restore the selection/active cell to match the refreshed export
# (first empty row right after the last data row).
$ws.Range("A44").Select()
